$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.220.86"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -1.85%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.581.74"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.22%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.64"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.86%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  -0.32%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  +0.30%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.805.59"
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.585.34"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.82%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.24%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -1.28%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.66"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -0.51%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.226.60"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -1.72%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -1.12%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -0.25%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.78"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -3.46%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.84"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  -1.09%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -1.23%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.29"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +8.58%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.287.28"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  -0.38%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.603"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -1.70%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("E41").Value = "  +2.33%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.770"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  -2.76%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.56"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.69%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.717.52"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.27%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.77"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.02%  "
